# Generate Report for Archive
#
# The localization report is being regenerated: the shared "Status" string
# moves from "Ready for handoff" to "In Translation". That shared string is
# referenced on the Overview sheet (zh-cn/de-de status columns) as well as
# on each per-locale sheet's Status column, so updating the text once
# flows through to every cell that uses it.
#
# Regenerating the report also re-runs the column auto-fit for the Status
# columns (since the display text length changed), so their widths are
# refreshed to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status text everywhere it appears -------------------------
# Overview sheet: E2 (zh-cn status) and F2 (de-de status)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Per-locale sheets: C2 is the "Status" column
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- Re-fit the Status columns to the new (shorter) text -------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn")
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de")
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.5   # column C ("Status")
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.5   # column C ("Status")

Write-Output "Updated status text and refreshed column widths for Overview/zh-cn/de-de sheets."
